$wb = $excel.ActiveWorkbook

# Insert the new "SCMCreds" worksheet right before the "Table" sheet
# (i.e. as the 4th sheet, after "UserPageData").
$tableSheet = $wb.Worksheets.Item("Table")
$newSheet = $wb.Worksheets.Add($tableSheet)
$newSheet.Name = "SCMCreds"

# Populate the data for the new sheet.
$newSheet.Range("A1").Value = "Mr"
$newSheet.Range("A2").Value = 1234567890
$newSheet.Range("A3").Value = "PO 45445, NY, 27756"
$newSheet.Range("A4").Value = 34

# Match the look of the other data sheets: size column A to its contents.
[void]$newSheet.Columns("A:A").AutoFit()

# Leave the selection just below the data, and make the new sheet active.
[void]$newSheet.Range("A5").Select()
[void]$newSheet.Select()
